$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple value updates (column B, and one column C cell)
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 17
$ws.Range("B3").Value = 47581786
$ws.Range("B4").Value = 76591981750000000
$ws.Range("B5").Value = 98410.14599999999
$ws.Range("B6").Value = 112.73258
$ws.Range("B7").Value = 88.67249200000001
$ws.Range("B8").Value = 2.824631
$ws.Range("B9").Value = 161.83235
$ws.Range("B10").Value = 92.32836
$ws.Range("B11").Value = 605.33724
$ws.Range("B12").Value = 2073.373
$ws.Range("B13").Value = 674.1830200000001
$ws.Range("B14").Value = 12001323.8
$ws.Range("B15").Value = 10937070
$ws.Range("B16").Value = 341142.92
$ws.Range("B17").Value = 806677.05
$ws.Range("B18").Value = 926.81639
$ws.Range("B19").Value = 362.91416
$ws.Range("B20").Value = 671618.5699999999
$ws.Range("B21").Value = 2568815

# Row 22: the formatting (font color / style) of B22 and C22 are swapped,
# and B22's value is also updated. Use a scratch cell to swap the styles
# via copy/paste-special (formats only), then clean up the scratch cell.
$b22 = $ws.Range("B22")
$c22 = $ws.Range("C22")
$scratch = $ws.Range("Z1")

$b22.Copy()
$scratch.PasteSpecial(-4122)

$c22.Copy()
$b22.PasteSpecial(-4122)

$scratch.Copy()
$c22.PasteSpecial(-4122)

$scratch.Clear()

$ws.Range("B22").Value = 238.2233

$ws.Range("B23").Value = 235.40831
$ws.Range("B24").Value = 137.92548
$ws.Range("B25").Value = 431.40883
$ws.Range("B26").Value = 487.85666
$ws.Range("B27").Value = 412.19445
$ws.Range("B28").Value = 435.87083
$ws.Range("B29").Value = 529.1650100000001
$ws.Range("B30").Value = 455.38336
$ws.Range("B31").Value = 596.36086
$ws.Range("B32").Value = 1430768.39
